$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "předmět" (subject) column values to abbreviated forms
$ws.Range("D2").Value = "IN"
$ws.Range("D3").Value = "PG"
$ws.Range("D4").Value = "IN"

# Update "zkoušející" (examiner) column values to last-name only
$ws.Range("F2").Value = "Bajer"
$ws.Range("F3").Value = "Šilar"
$ws.Range("F4").Value = "Štěpánek"

# Update "známka" (grade) for row 4 from text "N" to numeric 3
$ws.Range("E4").Value = 3
